$wb = $excel.ActiveWorkbook

# The "修改时间" (last-modified) column holds values like "202509211531" that
# are stored as TEXT (not numbers) in the workbook. A plain
# `.Value = "202509211537"` assignment would let Excel auto-detect the
# all-digit string and silently store it as a Number, which both changes the
# cell's type and (since the number is large) its displayed form. To keep
# the cell genuinely text we briefly force a Text number format, assign the
# value, then clear the format again so no extra cell styling lingers.
$newTimestamp = "202509211537"

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Sheet 1: 大智投资组合 - timestamps in column E, rows 2-9
$ws1 = $wb.Worksheets.Item("大智投资组合")
for ($r = 2; $r -le 9; $r++) {
    Set-TextValue $ws1.Cells.Item($r, 5) $newTimestamp
}

# Sheet 2: 大成投资组合 - timestamps in column E, rows 2-11
$ws2 = $wb.Worksheets.Item("大成投资组合")
for ($r = 2; $r -le 11; $r++) {
    Set-TextValue $ws2.Cells.Item($r, 5) $newTimestamp
}

# Sheet 3: 我的投资组合 - timestamps in column G, rows 2-13
$ws3 = $wb.Worksheets.Item("我的投资组合")
for ($r = 2; $r -le 13; $r++) {
    Set-TextValue $ws3.Cells.Item($r, 7) $newTimestamp
}
